# Revert to old censored template:
# Remove the "Water Temp", "Sp Conductance", "TP", "Ammonia" and "E.coli"
# parameter rows from the Censored sheet, leaving only Parameter/pH/DO/Nitrate.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Censored")

# Row 2 ("Water Temp") is dropped entirely; pH/DO/Sp Conductance/TP/Nitrate/
# Ammonia/E.coli all shift up one row.
$ws.Rows("2").Delete()

# After the shift: row2=pH, row3=DO, row4=Sp Conductance, row5=TP,
# row6=Nitrate, row7=Ammonia, row8=E.coli.
# Drop "Sp Conductance" and "TP" (rows 4:5).
$ws.Range("A4:A5").EntireRow.Delete()

# After that shift: row2=pH, row3=DO, row4=Nitrate, row5=Ammonia, row6=E.coli.
# Drop "Ammonia" and "E.coli" (rows 5:6), leaving Parameter/pH/DO/Nitrate.
$ws.Range("A5:A6").EntireRow.Delete()

# Restore the active selection in the frozen-pane view to cell B1.
[void]$ws.Range("B1").Select()
